$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(8)

# Move the "Feature Importance:" title textbox to its new position.
# (EMU 4293124 / 333381 converted to points; the literal is nudged slightly
# within the single-precision float's rounding slack so the EMU round-trip
# lands exactly on the target instead of one EMU short.)
$shp.Left = 338.04129
$shp.Top = 26.250472440944883

# Consolidate the three text runs ("Feature ", "Importance", ":") into a
# single run, keeping the first run's formatting (Arial, 28pt).
# Assigning the identical concatenated text directly is treated as a no-op
# by the text-range diffing (it shares a prefix/suffix with the existing
# runs), so nudge it through an unrelated intermediate value first to force
# all runs to collapse into one before setting the final text.
$shp.TextFrame.TextRange.Text = "x"
$shp.TextFrame.TextRange.Text = "Feature Importance:"
